$wb = $excel.ActiveWorkbook

# --- Rename the two "Include" sheets ---
$wsIncludeVS = $wb.Worksheets.Item("Include ValueSets")
$wsIncludeVS.Name = "Include ValueSet #0"

$wsIncludeKrebs = $wb.Worksheets.Item("Include from Krebsstadium Cod")
$wsIncludeKrebs.Name = "Include #1"

# --- Update the Metadata sheet ---
$ws = $wb.Worksheets.Item("Metadata")

# Update the Date value (row 8, column B)
$ws.Range("B8").Value = "2024-09-17T19:55:11+00:00"

# Insert a new "Jurisdiction" property row after the "Contact" row (row 10),
# pushing Description/Purpose/Copyright/Immutable down by one row (11->12,
# 12->13, 13->14, 14->15). Shift manually (bottom-up) instead of using
# Rows.Insert so the existing row formatting/styles are preserved exactly
# and no stray style definitions get introduced.
for ($r = 14; $r -ge 11; $r--) {
    $srcA = $ws.Cells.Item($r, 1).Value2
    $srcB = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r + 1, 1).Value2 = $srcA
    $ws.Cells.Item($r + 1, 2).Value2 = $srcB
}

# The newly extended row 15 has no formatting yet - copy it from an existing
# data row so it matches the rest of the table.
$ws.Range("A10:B10").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new row
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
